$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0.7354849458108593

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.7383857496430754

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = -0.7943852212231091
$ws.Range("D4").Value = -0.6858767186577983

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.6437900668177866
$ws.Range("D5").Value = -0.6523659492315754

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0.7531522638435943

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.8138762503080088

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = -0.7026441011123798
$ws.Range("D8").Value = -0.7150381880680573

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = -0.6507110192375765
$ws.Range("D9").Value = 0.6951854922468259
